$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "systemMenu" row -- the new menu button added to the game
$ws.Cells.Item(36, 1).Value = "systemMenu"

# Column B carries the Chinese display label and reuses the same
# "Chinese font" formatting already used by the other label cells
# (e.g. B35) -- copy that formatting over before writing the value.
$ws.Cells.Item(35, 2).Copy() | Out-Null
$ws.Cells.Item(36, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(36, 2).Value = "主界面"

$ws.Cells.Item(36, 3).Value = "lab_systemMenu"

# Leave the selection where the author ended up after entering the data
$ws.Range("B32").Select() | Out-Null
